$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 15: continuation of the calendar table (day 30, Sunday), with a time
# entry in C15 (same h:mm format already used by C14/D14) and "HK" tag in E15.
$ws.Range("A15").Value = 30
$ws.Range("B15").Value = "Sun"
$ws.Range("C15").Value = 0.20833333333333334
$ws.Range("C15").NumberFormat = "h:mm"
$ws.Range("E15").Value = "HK"

# Row 19: header row for a small expenses table. Row 20: matching data row.
$ws.Range("B19").Value = "Citi"
$ws.Range("A20").Value = "15-days"
$ws.Range("A19").Value = "Insur/Px"
$ws.Range("C19").Value = "Zurich"
$ws.Range("D19").Value = "Bluecross"

$ws.Range("B20").Value = 264.10000000000002
$ws.Range("C20").Formula = "=345*0.65"
$ws.Range("D20").Formula = "=335*0.65"

$ws.Range("A20").Select()
